# Generate Report for Handback
# Refresh the handoff/handback datetime stamps for the "8fc7b040..." file row
# (row 2) on both the zh-cn and de-de report sheets, as produced by a new
# CI run of the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 12:35:53"
$wsZhCn.Range("H2").Value = "2016-03-19 12:36:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 12:35:56"
$wsDeDe.Range("H2").Value = "2016-03-19 12:36:18"
